$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "C:\Users\Veeraraju_elluru\Desktop\Veeraraju\Personal\ATREE\test_images\download (1).png"
$ws.Range("B4").Value = "open3d,library"
$ws.Range("C4").Value = "picture"
$ws.Range("D4").Value = "None"
